$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Table on slide 16 (the "C1, C2 and C3" table) gets switched from
#    the deck's default table style to a different built-in table
#    style (identified by its style GUID).
# ------------------------------------------------------------------
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table
$tbl.ApplyStyle("{845C5E0D-CA58-46C8-A595-E57D0E0E3889}")

# ------------------------------------------------------------------
# 2) The presentation's design theme is switched from the custom
#    "Integral" theme to the built-in default "Office Theme" colour
#    scheme (Design tab -> Office Theme). Apply the Office Theme's
#    standard 12 theme colours to the deck's colour scheme.
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$cs = $slide1.ThemeColorScheme
$cs.Colors(1).RGB  = 0        # Dark 1      #000000
$cs.Colors(2).RGB  = 16777215 # Light 1     #FFFFFF
$cs.Colors(3).RGB  = 6968388  # Dark 2      #44546A
$cs.Colors(4).RGB  = 15132391 # Light 2     #E7E6E6
$cs.Colors(5).RGB  = 13998939 # Accent 1    #5B9BD5
$cs.Colors(6).RGB  = 3243501  # Accent 2    #ED7D31
$cs.Colors(7).RGB  = 10855845 # Accent 3    #A5A5A5
$cs.Colors(8).RGB  = 49407    # Accent 4    #FFC000
$cs.Colors(9).RGB  = 12874308 # Accent 5    #4472C4
$cs.Colors(10).RGB = 4697456  # Accent 6    #70AD47
$cs.Colors(11).RGB = 12673797 # Hyperlink   #0563C1
$cs.Colors(12).RGB = 7491477  # Followed Hyperlink #954F72
